$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add the new one after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "DevProfile"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "UserPersona"

# --- DevProfile sheet: insert a new column E (IsActive), shifting old E (Description) to F ---
$ws1.Columns.Item(5).Insert() | Out-Null

$ws1.Range("E1").Value = "IsActive"
$ws1.Range("E2").Value = $true
$ws1.Range("E3").Value = $true
$ws1.Range("E4").Value = $true
$ws1.Range("E5").Value = $true

# --- New rows 6 and 7 ---
$ws1.Range("A6").Value = "AMBEXRM01"
$ws1.Range("B6").Value = "Ambulatory Exam Room"
$ws1.Range("C6").Value = $true
$ws1.Range("D6").Value = $true
$ws1.Range("E6").Value = $true

$ws1.Range("A7").Value = "ROV01"
$ws1.Range("B7").Value = " Epic Rover "
$ws1.Range("C7").Value = $false
$ws1.Range("D7").Value = $false
$ws1.Range("E7").Value = $true
$ws1.Range("F7").Value = "Includes"

# --- New description values for rows 4 and 5 (column F) ---
$ws1.Range("F4").Value = "Includes Thin Client, 24`" Monitor with Keyboard, Mouse, proximity card reader, biometric scanner"
$ws1.Range("F5").Value = "Includes PC, 24`" Monitor with Keyboard, Mouse, Esignature Pad"

# --- Column widths: col E should match B-D (27.1640625), col F keeps 36.83203125 ---
$ws1.Columns.Item(5).ColumnWidth = 26.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 36.0

# --- Selection / view state ---
$ws1.Range("F6").Select() | Out-Null

# --- UserPersona sheet content ---
$ws2.Range("A1").Value = "UserPersonaCode"
$ws2.Range("B1").Value = "UserPersona"
$ws2.Range("C1").Value = "Associated Roles"

$ws2.Columns.Item(1).ColumnWidth = 17.833333333333332
$ws2.Columns.Item(3).ColumnWidth = 22.166666666666668

$ws2.Range("A2").Select() | Out-Null
$ws2.Activate() | Out-Null


